$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.175.34"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.657.44"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'215.43"
$ws.Range("D6").Value = "'0.5240"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'0.2629"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").Value = "'0.07770"
$ws.Range("E11").Value = "  +3.05%  "
$ws.Range("D12").Value = "1.662.68"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "'4.463"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "1.882.68"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "'0.5527"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").Value = "0.0₅8276"
$ws.Range("E16").Value = "  +4.13%  "
$ws.Range("D17").Value = "'65.18"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "26.177.62"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "'4.760"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "'190.51"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("D23").Value = "'6.373"
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'143.15"
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("D26").Value = "'0.1253"
$ws.Range("E26").Value = "  +3.18%  "
$ws.Range("D27").Value = "'7.416"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "'16.00"
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("D29").Value = "'1.430"
$ws.Range("E29").Value = "  +2.74%  "
$ws.Range("D30").Value = "'0.06117"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").Value = "'3.528"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").Value = "'3.423"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").Value = "'0.9998"
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").Value = "'2.401"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").Value = "'2.760"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'0.5663"
$ws.Range("E38").Value = "  -3.72%  "
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("D40").Value = "'5.912"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").Value = "'0.8552"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").Value = "1.032.11"
$ws.Range("E43").Value = "  -6.49%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "1.804.81"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").Value = "'56.12"
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("D48").Value = "'1.003"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").Value = "'8.073"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").Value = "'0.05164"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").Value = "'5.981"
$ws.Range("E51").Value = "  +2.20%  "

Write-Output "Updated cryptos list"
